# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '42.925.22'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '2.399.22'
$ws.Range('E3').Value = '  +4.72%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = '''335.18'
$ws.Range('E5').Value = '  +8.68%  '
$ws.Range('D6').Value = '''100.57'
$ws.Range('E6').Value = '  -10.42%  '
$ws.Range('D7').Value = '''0.644'
$ws.Range('E7').Value = '  +1.63%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '''0.637'
$ws.Range('E9').Value = '  +3.49%  '
$ws.Range('D10').Value = '''40.76'
$ws.Range('E10').Value = '  -8.10%  '
$ws.Range('D11').Value = '''0.0933'
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('D12').Value = '''8.54'
$ws.Range('E12').Value = '  -3.18%  '
$ws.Range('D13').Value = '''1.04'
$ws.Range('E13').Value = '  -3.77%  '
$ws.Range('D14').Value = '''16.94'
$ws.Range('E14').Value = '  +9.07%  '
$ws.Range('D15').Value = '''0.106'
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('D16').Value = '2.762.26'
$ws.Range('E16').Value = '  +4.54%  '
$ws.Range('D17').Value = '2.400.94'
$ws.Range('E17').Value = '  +4.95%  '
$ws.Range('D18').Value = '42.923.39'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = '''7.65'
$ws.Range('E19').Value = '  +6.21%  '
$ws.Range('D20').Value = '''0.0000107'
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('D21').Value = '''3.90'
$ws.Range('E21').Value = '  +10.47%  '
$ws.Range('D22').Value = '''76.41'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = '''269.06'
$ws.Range('E23').Value = '  +4.84%  '
$ws.Range('D24').Value = '''2.39'
$ws.Range('E24').Value = '  -2.95%  '
$ws.Range('D25').Value = '''10.21'
$ws.Range('E25').Value = '  +13.36%  '
$ws.Range('D26').Value = '''11.81'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('D28').Value = '''24.22'
$ws.Range('E28').Value = '  +8.15%  '
$ws.Range('E29').Value = '  -1.36%  '
$ws.Range('D30').Value = '''174.36'
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('D31').Value = '''3.12'
$ws.Range('E31').Value = '  -1.96%  '
$ws.Range('D32').Value = '''0.0920'
$ws.Range('E32').Value = '  +2.09%  '
$ws.Range('D33').Value = '''36.15'
$ws.Range('E33').Value = '  -7.41%  '
$ws.Range('D34').Value = '''5.98'
$ws.Range('E34').Value = '  +4.28%  '
$ws.Range('D35').Value = '''0.134'
$ws.Range('E35').Value = '  +3.34%  '
$ws.Range('D36').Value = '''4.77'
$ws.Range('E36').Value = '  -6.17%  '
$ws.Range('D37').Value = '''0.0362'
$ws.Range('E37').Value = '  -3.77%  '
$ws.Range('D38').Value = '''3.93'
$ws.Range('E38').Value = '  -5.52%  '
$ws.Range('E39').Value = '  +3.32%  '
$ws.Range('D40').Value = '''2.89'
$ws.Range('E40').Value = '  +12.75%  '
$ws.Range('E41').Value = '  +8.00%  '
$ws.Range('D42').Value = '''0.233'
$ws.Range('E42').Value = '  +0.59%  '
$ws.Range('D43').Value = '''69.75'
$ws.Range('E43').Value = '  -3.40%  '
$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D44').Value = '''94.51'
$ws.Range('E44').Value = '  +47.49%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '''1.00'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').Value = '''119.17'
$ws.Range('E46').Value = '  +9.89%  '
$ws.Range('D47').Value = '''11.98'
$ws.Range('E47').Value = '  -3.64%  '
$ws.Range('D48').Value = '''5.53'
$ws.Range('E48').Value = '  -2.80%  '
$ws.Range('D49').Value = '''9.18'
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('D50').Value = '1.644.44'
$ws.Range('E50').Value = '  +11.12%  '
$ws.Range('E51').Value = '  -1.21%  '
